$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 6050.6924
$ws.Range("I64").Value = 4870.375
$ws.Range("J64").Value = 7939.2
$ws.Range("K64").Value = 4870.375
$ws.Range("L64").Value = 7939.2
$ws.Range("M64").Value = -4622.375
$ws.Range("N64").Value = -8435.200000000001
$ws.Range("H67").Value = 6050.6924
$ws.Range("I67").Value = 4870.375
$ws.Range("J67").Value = 7939.2
$ws.Range("K67").Value = 4870.375
$ws.Range("L67").Value = 7939.2
$ws.Range("M67").Value = -4012.375
$ws.Range("N67").Value = -9655.200000000001
$ws.Range("H96").Value = 1452486.2
$ws.Range("I96").Value = 2202.25
$ws.Range("J96").Value = 2419342.2
$ws.Range("K96").Value = 6606.75
$ws.Range("L96").Value = 7258026.600000001
$ws.Range("M96").Value = -5233.75
$ws.Range("N96").Value = -7260772.600000001
$ws.Range("H125").Value = 6274.3335
$ws.Range("I125").Value = 679.75
$ws.Range("J125").Value = 10750
$ws.Range("K125").Value = 6117.75
$ws.Range("L125").Value = 96750
$ws.Range("M125").Value = -3657.75
$ws.Range("N125").Value = -101670
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5210.68
$ws.Range("I32").Value = 5114.6743
$ws.Range("K32").Value = 5114.6743
$ws.Range("M32").Value = -4827.6743
$ws.Range("H110").Value = 5549.263
$ws.Range("I110").Value = 5958.4614
$ws.Range("J110").Value = 4662.6665
$ws.Range("K110").Value = 5958.4614
$ws.Range("L110").Value = 4662.6665
$ws.Range("M110").Value = -3913.4614
$ws.Range("N110").Value = -8752.666499999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 99999
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H86").Value = 40615.723
$ws.Range("I86").Value = 57890
$ws.Range("J86").Value = 6067.1665
$ws.Range("K86").Value = 57890
$ws.Range("L86").Value = 6067.1665
$ws.Range("M86").Value = -56767
$ws.Range("N86").Value = -8313.166499999999
$ws.Range("H89").Value = 40615.723
$ws.Range("I89").Value = 57890
$ws.Range("J89").Value = 6067.1665
$ws.Range("K89").Value = 289450
$ws.Range("L89").Value = 30335.8325
$ws.Range("M89").Value = -283834
$ws.Range("N89").Value = -41567.8325
$ws.Range("H109").Value = 99998.5
$ws.Range("J109").Value = 99998.5
$ws.Range("L109").Value = 99998.5
$ws.Range("N109").Value = -102772.5
$ws.Range("H134").Value = 4349361
$ws.Range("I134").Value = 1604.8636
$ws.Range("K134").Value = 4814.5908
$ws.Range("M134").Value = -2279.5908
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 1046.6666
$ws.Range("I3").Value = 1046.6666
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1046.6666
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -933.6666
$ws.Range("N3").ClearContents()
$ws.Range("H4").Value = 4980
$ws.Range("I4").Value = 4900
$ws.Range("K4").Value = 4900
$ws.Range("M4").Value = -4788
$ws.Range("H7").Value = 31.4
$ws.Range("I7").Value = 28.8
$ws.Range("J7").Value = 32.7
$ws.Range("K7").Value = 28.8
$ws.Range("L7").Value = 32.7
$ws.Range("M7").Value = 84.2
$ws.Range("N7").Value = -258.7
$ws.Range("H31").Value = 23260220
$ws.Range("I31").Value = 47621624
$ws.Range("J31").Value = 6149.9546
$ws.Range("K31").Value = 47621624
$ws.Range("L31").Value = 6149.9546
$ws.Range("M31").Value = -47621329
$ws.Range("N31").Value = -6739.9546
$ws.Range("H34").Value = 23260220
$ws.Range("I34").Value = 47621624
$ws.Range("J34").Value = 6149.9546
$ws.Range("K34").Value = 47621624
$ws.Range("L34").Value = 6149.9546
$ws.Range("M34").Value = -47621422
$ws.Range("N34").Value = -6553.9546
$ws.Range("H99").Value = 14412.521
$ws.Range("I99").Value = 8421.611000000001
$ws.Range("K99").Value = 8421.611000000001
$ws.Range("M99").Value = -6923.611000000001
$ws.Range("H126").Value = 14412.521
$ws.Range("I126").Value = 8421.611000000001
$ws.Range("K126").Value = 25264.833
$ws.Range("M126").Value = -22794.833
$ws.Range("H132").Value = 1707.5151
$ws.Range("I132").Value = 1707.5151
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5122.5453
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2592.5453
$ws.Range("N132").ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 542.3333
$ws.Range("I5").Value = 515.8
$ws.Range("J5").Value = 675
$ws.Range("K5").Value = 1547.4
$ws.Range("L5").Value = 2025
$ws.Range("M5").Value = -1435.4
$ws.Range("N5").Value = -2249
$ws.Range("H12").Value = 2400.3635
$ws.Range("I12").Value = 1334.4445
$ws.Range("J12").Value = 3138.3076
$ws.Range("K12").Value = 4003.3335
$ws.Range("L12").Value = 9414.9228
$ws.Range("M12").Value = -3830.3335
$ws.Range("N12").Value = -9760.9228
$ws.Range("H14").Value = 13889.917
$ws.Range("I14").Value = 13889.917
$ws.Range("K14").Value = 41669.751
$ws.Range("M14").Value = -41496.751
$ws.Range("H81").Value = 30411
$ws.Range("J81").Value = 30411
$ws.Range("L81").Value = 91233
$ws.Range("N81").Value = -93479
$ws.Range("H84").Value = 30411
$ws.Range("J84").Value = 30411
$ws.Range("L84").Value = 273699
$ws.Range("N84").Value = -284931
$ws.Range("H86").Value = 392
$ws.Range("I86").Value = 385
$ws.Range("J86").Value = 402.5
$ws.Range("K86").Value = 1155
$ws.Range("L86").Value = 1207.5
$ws.Range("M86").Value = 31
$ws.Range("N86").Value = -3579.5
$ws.Range("H89").Value = 392
$ws.Range("I89").Value = 385
$ws.Range("J89").Value = 402.5
$ws.Range("K89").Value = 3465
$ws.Range("L89").Value = 3622.5
$ws.Range("M89").Value = 2463
$ws.Range("N89").Value = -15478.5
$ws.Range("H135").Value = 542.3333
$ws.Range("I135").Value = 515.8
$ws.Range("J135").Value = 675
$ws.Range("K135").Value = 4642.2
$ws.Range("L135").Value = 6075
$ws.Range("M135").Value = -2107.2
$ws.Range("N135").Value = -11145
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 47777
$ws.Range("J15").Value = 59999
$ws.Range("L15").Value = 59999
$ws.Range("N15").Value = -60575
$ws.Range("H81").Value = 47777
$ws.Range("J81").Value = 59999
$ws.Range("L81").Value = 59999
$ws.Range("N81").Value = -61995
$ws.Range("H84").Value = 47777
$ws.Range("J84").Value = 59999
$ws.Range("L84").Value = 179997
$ws.Range("N84").Value = -189981
$ws.Range("H126").Value = 3865.75
$ws.Range("I126").Value = 3194.5715
$ws.Range("J126").Value = 4805.4
$ws.Range("K126").Value = 9583.7145
$ws.Range("L126").Value = 14416.2
$ws.Range("M126").Value = -7113.7145
$ws.Range("N126").Value = -19356.2
$ws.Range("H136").Value = 9886
$ws.Range("J136").Value = 9886
$ws.Range("L136").Value = 29658
$ws.Range("N136").Value = -34758
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6452.1816
$ws.Range("I7").Value = 2996.7144
$ws.Range("J7").Value = 12499.25
$ws.Range("K7").Value = 2996.7144
$ws.Range("L7").Value = 12499.25
$ws.Range("M7").Value = -2884.7144
$ws.Range("N7").Value = -12723.25
$ws.Range("H8").Value = 200000
$ws.Range("J8").Value = 200000
$ws.Range("L8").Value = 200000
$ws.Range("N8").Value = -200280
$ws.Range("H80").Value = 72324.664
$ws.Range("J80").Value = 72324.664
$ws.Range("L80").Value = 72324.664
$ws.Range("N80").Value = -74570.664
$ws.Range("H83").Value = 72324.664
$ws.Range("J83").Value = 72324.664
$ws.Range("L83").Value = 216973.992
$ws.Range("N83").Value = -228205.992
$ws.Range("H93").Value = 1794436.8
$ws.Range("I93").Value = 894.5
$ws.Range("K93").Value = 894.5
$ws.Range("M93").Value = 353.5
$ws.Range("H126").Value = 6452.1816
$ws.Range("I126").Value = 2996.7144
$ws.Range("J126").Value = 12499.25
$ws.Range("K126").Value = 8990.143199999999
$ws.Range("L126").Value = 37497.75
$ws.Range("M126").Value = -6520.143199999999
$ws.Range("N126").Value = -42437.75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 99476
$ws.Range("J99").Value = 99476
$ws.Range("L99").Value = 99476
$ws.Range("N99").Value = -105466
$ws.Range("H122").Value = 2681.889
$ws.Range("I122").Value = 2571.5293
$ws.Range("K122").Value = 7714.5879
$ws.Range("M122").Value = -5264.5879
$ws.Range("H126").Value = 2002.6857
$ws.Range("I126").Value = 1665.5667
$ws.Range("K126").Value = 4996.7001
$ws.Range("M126").Value = -2526.7001
$ws.Range("H132").Value = 176279.53
$ws.Range("I132").Value = 809
$ws.Range("K132").Value = 2427
$ws.Range("M132").Value = 103
